# Applies the edit described by the diff:
#  1. Remove the standalone "Meta description: ..." paragraph that used to sit
#     right under the title.
#  2. Insert a new bold paragraph containing the page title
#     ("Play Cai Shen Dao slot free: Review and analysis") right before the
#     final (italic) paragraph.
#  3. Replace the text of that final italic paragraph (previously the
#     "Prompt: ..." image-generation prompt) with the meta-description text
#     that used to live near the top of the document.

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph -----------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Meta description")) {
        $p.Range.Delete()
        break
    }
}

# --- Step 2: insert a new bold title paragraph just before the final -------
#             italic "Prompt: ..." paragraph.
$promptIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith("Prompt:")) {
        $promptIndex = $i
        break
    }
}

$anchor = $d.Paragraphs($promptIndex - 1)  # paragraph right before "Prompt: ..."
$anchorRange = $anchor.Range
$anchorRange.Collapse(0)                   # collapse to the end of that paragraph

$insertStart = $anchorRange.End
$titleText = "Play Cai Shen Dao slot free: Review and analysis"

$anchorRange.InsertAfter([char]13)        # new paragraph mark

$titleRange = $d.Range($insertStart, $insertStart)
$titleRange.InsertAfter($titleText)

$boldRange = $d.Range($insertStart, $insertStart + $titleText.Length)
$boldRange.Font.Bold = $true

# --- Step 3: swap the old "Prompt: ..." text for the meta description text -
$oldPrompt = 'Prompt: Create a fun and engaging feature image for "Cai Shen Dao" that fits the description of a happy Maya warrior with glasses. The image should be in cartoon style and showcase the Chinese culture theme of the game, with symbols such as the carp, fan, and amulets. The image should also have a touch of humor and a playful vibe to appeal to the game''s target audience.'
$newMeta = 'Explore Cai Shen Dao slot game with our review. Play for free and discover its impressive graphics, thematic symbols, medium volatility, and free spin feature.'

$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newMeta, 2)
